# Week 13: release 0.2
# Adds "Question 7" and "Question 9" worked-example tables to the T8 sheet,
# widens column A, and updates the selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T8")

# --- Column A width: 24.28515625 -> 43.28515625 (character-width units) ---
$ws.Columns("A").ColumnWidth = 42.5

# ---------------------------------------------------------------------
# Question 7 block (rows 108-113)
# ---------------------------------------------------------------------
$ws.Range("A108").Value = "Question 7"
$ws.Range("I108").Value = "SUM"

$ws.Range("A109").Value = "Let 2004 = year 0, Year, X"
$ws.Range("B109").Value = 0
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 2
$ws.Range("E109").Value = 3
$ws.Range("F109").Value = 4
$ws.Range("G109").Value = 5
$ws.Range("H109").Value = 6
$ws.Range("I109").Formula = "=SUM(B109:H109)"

$ws.Range("A110").Value = "Average annual earnings (`$ '000), Y"
$ws.Range("B110").Value = 59
$ws.Range("C110").Value = 70
$ws.Range("D110").Value = 77
$ws.Range("E110").Value = 87
$ws.Range("F110").Value = 89
$ws.Range("G110").Value = 122
$ws.Range("H110").Value = 137
$ws.Range("I110").Formula = "=SUM(B110:H110)"

$ws.Range("A111").Value = "XY"
$ws.Range("B111").Formula = "=B109*B110"
$ws.Range("C111").Formula = "=C109*C110"
$ws.Range("D111").Formula = "=D109*D110"
$ws.Range("E111").Formula = "=E109*E110"
$ws.Range("F111").Formula = "=F109*F110"
$ws.Range("G111").Formula = "=G109*G110"
$ws.Range("H111").Formula = "=H109*H110"
$ws.Range("I111").Formula = "=SUM(B111:H111)"

$ws.Range("A112").Value = "X^2"
$ws.Range("B112").Formula = "=B109^2"
$ws.Range("C112").Formula = "=C109^2"
$ws.Range("D112").Formula = "=D109^2"
$ws.Range("E112").Formula = "=E109^2"
$ws.Range("F112").Formula = "=F109^2"
$ws.Range("G112").Formula = "=G109^2"
$ws.Range("H112").Formula = "=H109^2"
$ws.Range("I112").Formula = "=SUM(B112:H112)"

$ws.Range("A113").Value = "Y^2"
$ws.Range("B113").Formula = "=B110^2"
$ws.Range("C113").Formula = "=C110^2"
$ws.Range("D113").Formula = "=D110^2"
$ws.Range("E113").Formula = "=E110^2"
$ws.Range("F113").Formula = "=F110^2"
$ws.Range("G113").Formula = "=G110^2"
$ws.Range("H113").Formula = "=H110^2"
$ws.Range("I113").Formula = "=SUM(B113:H113)"

# Bold style for the label/sum columns, matching the rest of the sheet
$ws.Range("A108,I108,A109,I109,A110,I110,A111,I111,A112,I112,A113,I113").Font.Bold = $true

# ---------------------------------------------------------------------
# Question 9 block (rows 115-121)
# ---------------------------------------------------------------------
$ws.Range("A115").Value = "Question 9"

$ws.Range("A116").Value = "Let 31 = 0"
$ws.Range("L116").Value = "SUM"

$ws.Range("A117").Value = "Age (year), X"
$ws.Range("B117").Formula = "=31-31"
$ws.Range("C117").Formula = "=36-31"
$ws.Range("D117").Formula = "=39-31"
$ws.Range("E117").Formula = "=43-31"
$ws.Range("F117").Formula = "=47-31"
$ws.Range("G117").Formula = "=52-31"
$ws.Range("H117").Formula = "=58-31"
$ws.Range("I117").Formula = "=63-31"
$ws.Range("J117").Formula = "=69-31"
$ws.Range("K117").Formula = "=74-31"
$ws.Range("L117").Formula = "=SUM(B117:K117)"

$ws.Range("A118").Value = "Cholesterol Level, Y"
$ws.Range("B118").Value = 165
$ws.Range("C118").Value = 181
$ws.Range("D118").Value = 177
$ws.Range("E118").Value = 193
$ws.Range("F118").Value = 213
$ws.Range("G118").Value = 191
$ws.Range("H118").Value = 189
$ws.Range("I118").Value = 154
$ws.Range("J118").Value = 235
$ws.Range("K118").Value = 198
$ws.Range("L118").Formula = "=SUM(B118:K118)"

$ws.Range("A119").Value = "XY"
$ws.Range("B119").Formula = "=B117*B118"
$ws.Range("C119").Formula = "=C117*C118"
$ws.Range("D119").Formula = "=D117*D118"
$ws.Range("E119").Formula = "=E117*E118"
$ws.Range("F119").Formula = "=F117*F118"
$ws.Range("G119").Formula = "=G117*G118"
$ws.Range("H119").Formula = "=H117*H118"
$ws.Range("I119").Formula = "=I117*I118"
$ws.Range("J119").Formula = "=J117*J118"
$ws.Range("K119").Formula = "=K117*K118"
$ws.Range("L119").Formula = "=SUM(B119:K119)"

$ws.Range("A120").Value = "X^2"
$ws.Range("B120").Formula = "=B117^2"
$ws.Range("C120").Formula = "=C117^2"
$ws.Range("D120").Formula = "=D117^2"
$ws.Range("E120").Formula = "=E117^2"
$ws.Range("F120").Formula = "=F117^2"
$ws.Range("G120").Formula = "=G117^2"
$ws.Range("H120").Formula = "=H117^2"
$ws.Range("I120").Formula = "=I117^2"
$ws.Range("J120").Formula = "=J117^2"
$ws.Range("K120").Formula = "=K117^2"
$ws.Range("L120").Formula = "=SUM(B120:K120)"

$ws.Range("A121").Value = "Y^2"
$ws.Range("B121").Formula = "=B118^2"
$ws.Range("C121").Formula = "=C118^2"
$ws.Range("D121").Formula = "=D118^2"
$ws.Range("E121").Formula = "=E118^2"
$ws.Range("F121").Formula = "=F118^2"
$ws.Range("G121").Formula = "=G118^2"
$ws.Range("H121").Formula = "=H118^2"
$ws.Range("I121").Formula = "=I118^2"
$ws.Range("J121").Formula = "=J118^2"
$ws.Range("K121").Formula = "=K118^2"
$ws.Range("L121").Formula = "=SUM(B121:K121)"

$ws.Range("A115,A116,L116,L117,L118,A119,L119,A120,L120,A121,L121").Font.Bold = $true

# ---------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 108
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L121").Select()
